$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new notification row (row 9) to the table
$ws.Range("A9").Value = 9
$ws.Range("B9").Value = "Thông báo họp Dự án sáng 18/09/2022"
$ws.Range("C9").Value = "<p>Mời anh Dương Văn Hữu ....</p>"
$ws.Range("D9").Value = "Ban Đào Tạo"
$ws.Range("E9").Value = "18/09/2022 09:58"
$ws.Range("F9").Value = "https://us06web.zoom.us/postattendee?mn=bmYlpIgEcR-bI-lM1s-fjt1LNaCh6M98wrRh.9veeSYB9uE4OLuCg"
